$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels from the legacy UPPERCASE technical names to the
# new human-friendly CamelCase labels used by the XML export/import tool.
$ws.Range("B1").Value = "FirstName"
$ws.Range("C1").Value = "LastName"
$ws.Range("D1").Value = "DateofBirth"
$ws.Range("E1").Value = "ClasseName"
$ws.Range("F1").Value = "Phone"
$ws.Range("G1").Value = "Email"

# Switch the CNE numbering scheme to the new 8-digit format; the rest of
# the column (A3:A11) is driven by formulas relative to A2, so they pick
# up the new numbering automatically on recalculation.
$ws.Range("A2").Value = 20000001

# Leave the cursor where the author left it when they saved the file.
$ws.Range("C9").Select() | Out-Null
